$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the used A1:E51 range) used to stage values that look like
# plain numbers (e.g. "209.54") so they round-trip as TEXT, matching the source data
# (which stores every Price/Volume cell as an inline string), without leaving behind
# a changed NumberFormat/style on the real target cell.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = '27.681.40'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.619.26'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.07%  '
$helper.Value = '209.54'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("E6").Value = '  -0.85%  '
$helper.Value = '0.991'
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.06%  '
$helper.Value = '23.09'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("E9").Value = '  -0.86%  '
$helper.Value = '0.0606'
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -0.81%  '
$helper.Value = '0.0875'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '1.847.54'
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("D13").Value = '1.626.73'
$ws.Range("E14").Value = '  -1.17%  '
$helper.Value = '0.557'
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -0.70%  '
$helper.Value = '64.77'
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '27.687.06'
$ws.Range("E17").Value = '  -0.18%  '
$helper.Value = '227.52'
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").Value = '0.0₃0715'
$ws.Range("E20").Value = '  -0.81%  '
$helper.Value = '0.991'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.14%  '
$helper.Value = '4.31'
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -0.71%  '
$helper.Value = '10.03'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("E24").Value = '  -0.53%  '
$helper.Value = '154.83'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +0.49%  '
$helper.Value = '6.89'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '1.392.40'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("E35").Value = '  +1.63%  '
$helper.Value = '0.994'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("E40").Value = '  -2.59%  '
$helper.Value = '0.991'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  -1.44%  '
$ws.Range("E43").Value = '  -0.08%  '
$helper.Value = '65.52'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").Value = '1.756.40'
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("E47").Value = '  -3.37%  '
$helper.Value = '87.71'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("E50").Value = '  -0.58%  '
$helper.Value = '7.51'
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +0.92%  '

# Clean up the helper cell/clipboard state
$helper.Clear()
$excel.CutCopyMode = $false

